$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.844.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.083.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.38%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.085.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.768.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.55%  "
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.138"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.06%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("E32").Value = "  +4.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0629"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("E37").Value = "  +2.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0991"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.53%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.459.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.15%  "
$ws.Range("E48").Value = "  +4.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.276.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.65%  "
